# Applies the updated odds values to Sheet1 of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F2"  = 1.49
    "H2"  = 8.199999999999999
    "M2"  = 1.08
    "Q2"  = 2.18
    "F3"  = 3.8
    "G3"  = 4.5
    "I3"  = 2
    "G4"  = 4.6
    "I4"  = 1.92
    "J4"  = 3.85
    "N4"  = 4.3
    "P4"  = 2.14
    "R4"  = 1.45
    "X4"  = 21
    "Z4"  = 12.5
    "AA4" = 22
    "AB4" = 21
    "AF4" = 34
    "AK4" = 200
    "AN4" = 50
    "P5"  = 2.46
    "Q5"  = 1.63
    "R5"  = 1.6
    "S5"  = 2.56
    "U5"  = 2.7
    "X5"  = 24
    "Z5"  = 17.5
    "AA5" = 32
    "AE5" = 22
    "AI5" = 30
    "AN5" = 23
    "F6"  = 2.26
    "I6"  = 3.6
    "J6"  = 3.55
    "K6"  = 4.1
    "G7"  = 2.48
    "P8"  = 1.4
    "F10" = 2.02
    "G10" = 2.18
    "P10" = 1.42
    "Q10" = 2.96
    "G11" = 3.05
    "H11" = 3.2
    "I11" = 3.55
    "J11" = 2.68
    "K11" = 2.88
    "P11" = 1.35
    "Q11" = 3.5
    "G12" = 2.14
    "J12" = 3.1
    "Q12" = 2.1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
